$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "In Translation" -> "Ready for handoff"
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Handoff date/time text updates
$wsOverview.Range("D2").Value = "2016-29-19 16:29:57"
$wsZhCn.Range("E2").Value = "2016-03-19 16:29:54"
$wsDeDe.Range("E2").Value = "2016-03-19 16:29:57"
